# LCSC_Bom.xlsx update: "Updated BOMs to reflect changes"
#
# - Insert a new capacitor line (CL10A106MA8NRNC) after the CL10A105KB8NNNC row
# - Insert a new transistor line (BSS138 / KEXIN) after the MC-311D row
# - Bump the 0603WAF1002T5E quantity from 11 to 13
# - Re-source the AMS1117-3.3 regulator (manufacturer + LCSC part number)
# - Re-source the Seiko Epson crystal (manufacturer part number + LCSC part number)
# - Re-enter the Ext Quantity formula without the CEILING.MATH wrapper
# - Move the active selection to D28

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new BOM rows -----------------------------------------
# Row 9 (old) = CL10A105KB8NNNC -> new row goes in right after it.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 3).Value = "CL10A106MA8NRNC"
$ws.Cells.Item(9, 4).Value = "Samsung Electro-Mechanics"
$ws.Cells.Item(9, 6).Value = "C96446"

# Row 11 (after the above insert) = MC-311D -> new row goes in right after it.
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 3).Value = "BSS138"
$ws.Cells.Item(12, 4).Value = "KEXIN"
$ws.Cells.Item(12, 6).Value = "C489349"

# --- Quantity change: 0603WAF1002T5E 11 -> 13 -----------------------------
$ws.Cells.Item(14, 1).Value = 13

# --- Re-source AMS1117-3.3 --------------------------------------------
$ws.Cells.Item(21, 4).Value = "Advanced Monolithic Systems"
$ws.Cells.Item(21, 6).Value = "C6186"

# --- Re-source the Seiko Epson crystal ------------------------------------
$ws.Cells.Item(25, 3).Value = "Q13FC1350000200"
$ws.Cells.Item(25, 6).Value = "C48615"

# --- Re-enter the Ext Quantity (column B) formula on every data row -------
# (drops the CEILING.MATH wrapper and resets the cell style to Normal, same
# as what happens when Excel re-types/re-fills a formula over the old one)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Formula = "=(A$r*5)+(0.4*(A$r*5))"
    $cell.Style = "Normal"
}

# --- Move the selection -----------------------------------------------
$ws.Range("D28").Select()
